# Update gh-pages to output generated at 456a3b4
#
# Applies to both the "展览" (sheet1) and "全部类型" (sheet4) worksheets,
# which carried identical data:
#   1) Refresh the "want to go" counts (column F) for a handful of
#      existing events.
#   2) A brand-new event ("南昌·代号鸢盛花行only") is inserted as row 23;
#      the former row 23 ("九江·第三届ACD动漫游戏嘉年华") is pushed down to
#      row 24, with its running index (column A) bumped from 22 to 23 and
#      its "want to go" count refreshed from 7 to 14.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- 1) Plain "want to go" (column F) refreshes -----------------------
    $ws.Cells.Item(2, 6).Value = 1074
    $ws.Cells.Item(5, 6).Value = 3072
    $ws.Cells.Item(7, 6).Value = 2343
    $ws.Cells.Item(11, 6).Value = 1139
    $ws.Cells.Item(12, 6).Value = 35
    $ws.Cells.Item(13, 6).Value = 50
    $ws.Cells.Item(15, 6).Value = 808
    $ws.Cells.Item(16, 6).Value = 280
    $ws.Cells.Item(17, 6).Value = 301
    $ws.Cells.Item(19, 6).Value = 13
    $ws.Cells.Item(21, 6).Value = 54

    # --- 2) Push the old row 23 down to a new row 24 -----------------------
    # Copy column A's format (bold/centered/bordered style) down first so
    # the new row's index cell matches the look of every other row.
    $ws.Cells.Item(23, 1).Copy()
    $ws.Cells.Item(24, 1).PasteSpecial(-4122)

    $ws.Cells.Item(24, 1).Value = 23

    $ws.Cells.Item(24, 2).NumberFormat = "@"
    $ws.Cells.Item(24, 2).Value = "2024-05-01"
    $ws.Cells.Item(24, 2).ClearFormats()

    $ws.Cells.Item(24, 3).Value = "九江·第三届ACD动漫游戏嘉年华"
    $ws.Cells.Item(24, 4).Value = "九瑞大道与重庆路交汇处西南角 九江国际会展中心"

    $ws.Cells.Item(24, 5).NumberFormat = "@"
    $ws.Cells.Item(24, 5).Value = "2024.05.01 09:00-05.02 17:00"
    $ws.Cells.Item(24, 5).ClearFormats()

    $ws.Cells.Item(24, 6).Value = 14
    $ws.Cells.Item(24, 7).Value = 39.9
    $ws.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82464"
    $ws.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg"

    # --- 3) Overwrite row 23 in place with the new event -------------------
    # Column A (the running index, 22) is left untouched on purpose.
    $ws.Cells.Item(23, 2).NumberFormat = "@"
    $ws.Cells.Item(23, 2).Value = "2024-04-20"
    $ws.Cells.Item(23, 2).ClearFormats()

    $ws.Cells.Item(23, 3).Value = "南昌·代号鸢盛花行only"
    $ws.Cells.Item(23, 4).Value = "民德路411号 东方豪景花园酒店(民德路店)"

    $ws.Cells.Item(23, 5).NumberFormat = "@"
    $ws.Cells.Item(23, 5).Value = "2024.04.20 09:30-04.20 17:30"
    $ws.Cells.Item(23, 5).ClearFormats()

    $ws.Cells.Item(23, 6).Value = 1
    $ws.Cells.Item(23, 7).Value = 78
    $ws.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82529"
    $ws.Cells.Item(23, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png"
}
